# Updates cryptos list data (prices / 1h volume %, and two swapped
# coin-ranking rows) to match the refreshed scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "63.869.13"
$ws.Range("E2").Value = "  -0.02%  "

# Row 3
$ws.Range("D3").Value = "2.735.58"
$ws.Range("E3").Value = "  -0.56%  "

# Row 4
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "565.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.65%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "160.86"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.98%  "

# Row 8
$ws.Range("E8").Value = "  -0.96%  "

# Row 9
$ws.Range("E9").Value = "  -0.13%  "

# Row 10
$ws.Range("E10").Value = "  +4.25%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.63"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.53%  "

# Row 12
$ws.Range("E12").Value = "  -1.40%  "

# Row 13
$ws.Range("D13").Value = "3.221.49"
$ws.Range("E13").Value = "  -0.51%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.97"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.63%  "

# Row 15
$ws.Range("D15").Value = "63.704.25"
$ws.Range("E15").Value = "  +0.27%  "

# Row 16
$ws.Range("E16").Value = "  -0.74%  "

# Row 17
$ws.Range("D17").Value = "2.739.94"
$ws.Range("E17").Value = "  -0.50%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.32"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.86%  "

# Row 19
$ws.Range("E19").Value = "  -1.69%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "355.98"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.16%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.61"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.29%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.998"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.04%  "

# Row 23
$ws.Range("E23").Value = "  -2.56%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.18"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.50%  "

# Row 25
$ws.Range("E25").Value = "  -0.29%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.06%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.37"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.21%  "

# Row 28
$ws.Range("D28").Value = "0.0₃0910"
$ws.Range("E28").Value = "  -0.07%  "

# Row 29
$ws.Range("B29").Value = "Fetch.AI"
$ws.Range("C29").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.39"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +11.47%  "

# Row 30
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.61%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.17"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.34%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "167.35"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.54%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.92"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.14%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.49"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.05%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "20.06"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.51%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.81"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.09%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.979"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.85%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "348.42"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.24%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.30"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.22%  "

# Row 41
$ws.Range("E41").Value = "  -1.81%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "38.64"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.62%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.86"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.58%  "

# Row 44
$ws.Range("B44").Value = "Hedera"
$ws.Range("C44").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0584"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.59%  "

# Row 45
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "20.94"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.50%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.631"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.83%  "

# Row 47
$ws.Range("E47").Value = "  -0.74%  "

# Row 48
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "132.61"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.57%  "

# Row 49
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0995"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.33%  "

# Row 50
$ws.Range("E50").Value = "  -0.14%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.07"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.42%  "
